$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.577.77"
$ws.Range("E2").Value = "  -0.50%  "

$ws.Range("D3").Value = "1.623.80"
$ws.Range("E3").Value = "  -1.33%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.69"
$ws.Range("E5").Value = "  -0.74%  "

$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.31"
$ws.Range("E8").Value = "  +0.14%  "

$ws.Range("E9").Value = "  +2.32%  "

$ws.Range("E10").Value = "  -0.20%  "

$ws.Range("E11").Value = "  -0.37%  "

$ws.Range("D12").Value = "1.852.71"
$ws.Range("E12").Value = "  -1.42%  "

$ws.Range("D13").Value = "1.623.76"
$ws.Range("E13").Value = "  -1.43%  "

$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.551"
$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.26"
$ws.Range("E16").Value = "  +0.99%  "

$ws.Range("D17").Value = "27.529.11"
$ws.Range("E17").Value = "  -0.62%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.55"
$ws.Range("E18").Value = "  +0.46%  "

$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.58"
$ws.Range("E20").Value = "  -0.47%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.42"
$ws.Range("E22").Value = "  +3.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.34"
$ws.Range("E23").Value = "  +1.26%  "

$ws.Range("E24").Value = "  +6.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.22"
$ws.Range("E25").Value = "  +0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.89"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("E27").Value = "  -0.29%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.57"
$ws.Range("E28").Value = "  -0.35%  "

$ws.Range("E29").Value = "  -0.17%  "

$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("E32").Value = "  -0.43%  "

$ws.Range("D33").Value = "1.473.58"
$ws.Range("E33").Value = "  +2.26%  "

$ws.Range("E34").Value = "  -1.82%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.55"
$ws.Range("E35").Value = "  -2.54%  "

$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.955"
$ws.Range("E37").Value = "  +7.79%  "

$ws.Range("E38").Value = "  +0.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.873"
$ws.Range("E39").Value = "  -0.51%  "

$ws.Range("E40").Value = "  -1.91%  "

$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("E42").Value = "  -1.72%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "67.77"
$ws.Range("E43").Value = "  +1.10%  "

$ws.Range("E44").Value = "  -1.82%  "

$ws.Range("E45").Value = "  -1.72%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.28"
$ws.Range("E46").Value = "  -5.25%  "

$ws.Range("D47").Value = "1.764.36"
$ws.Range("E47").Value = "  -1.39%  "

$ws.Range("E48").Value = "  +1.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "87.44"
$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("E50").Value = "  -2.07%  "

$ws.Range("E51").Value = "  +1.92%  "
